$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.202.47"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.161.93"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("D9").Value = "3.160.02"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.509"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "3.679.71"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "66.254.49"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "3.158.88"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.726"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "503.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0418"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "0.0₃0672"
$ws.Range("E42").Value = "  +4.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.295"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "2.817.28"
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.14%  "
